$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 14 data rows (rows 2-15): line1..line6, extr1..extr8.
# The target sheet needs 16 data rows (rows 2-17): line1..line8, extr1..extr8.
# This is achieved by shifting the existing extr1..extr8 rows (old rows 8-15)
# down by two rows (to rows 10-17), then filling the freed rows 8-9 with the
# new line7 / line8 entries, and finally updating a handful of in_service
# flags that changed for some of the extr rows.

# Step 1: shift rows 8-15 down to rows 10-17 (process bottom-up so source
# rows are read before they get overwritten by the shift). Column A is a
# simple running index (row number - 2), so it is recomputed rather than
# copied verbatim.
for ($r = 15; $r -ge 8; $r--) {
    $destRow = $r + 2
    $ws.Cells.Item($destRow, 1).Value = $destRow - 2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($destRow, 5).Value = $ws.Cells.Item($r, 5).Value2
}

# Rows 16-17 are brand new territory (previously beyond the used range) -
# copy the formatting (bold/border/alignment) used by the rest of column A
# so they match the rest of the table.
$ws.Range("A7").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 2: populate the two freed rows with the new line7 / line8 data.
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Step 3: fix the in_service flags that changed for the shifted extr rows.
$ws.Range("E10").Value = $true   # extr1: in_service 0 -> 1
$ws.Range("E11").Value = $true   # extr2: in_service 0 -> 1
$ws.Range("E13").Value = $false  # extr4: in_service 1 -> 0
$ws.Range("E16").Value = $false  # extr7: in_service 1 -> 0

Write-Output "done"
